$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells remain plain text so values like
# "17.40" or "0.672" are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '34.823.09'
$ws.Range("E2").Value = '  -0.82%  '
$ws.Range("D3").Value = '1.838.01'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = '230.85'
$ws.Range("E5").Value = '  -0.71%  '
$ws.Range("E6").Value = '  +0.98%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '39.63'
$ws.Range("E8").Value = '  -3.92%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").Value = '0.0984'
$ws.Range("E11").Value = '  -1.50%  '
$ws.Range("D12").Value = '2.103.94'
$ws.Range("E12").Value = '  +1.21%  '
$ws.Range("D13").Value = '1.862.79'
$ws.Range("E13").Value = '  +2.54%  '
$ws.Range("D14").Value = '11.34'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").Value = '0.672'
$ws.Range("E15").Value = '  +1.88%  '
$ws.Range("D16").Value = '4.65'
$ws.Range("E16").Value = '  -0.09%  '
$ws.Range("D17").Value = '34.835.33'
$ws.Range("E17").Value = '  -0.74%  '
$ws.Range("D18").Value = '69.69'
$ws.Range("E18").Value = '  +0.09%  '
$ws.Range("E19").Value = '  -0.72%  '
$ws.Range("D20").Value = '240.68'
$ws.Range("E20").Value = '  +0.59%  '
$ws.Range("D21").Value = '12.17'
$ws.Range("E21").Value = '  +2.08%  '
$ws.Range("E22").Value = '  -0.19%  '
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").Value = '2.26'
$ws.Range("E24").Value = '  -0.65%  '
$ws.Range("D25").Value = '171.15'
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("E27").Value = '  +2.82%  '
$ws.Range("D28").Value = '17.40'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("E29").Value = '  -4.91%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  -0.26%  '
$ws.Range("D32").Value = '3.94'
$ws.Range("E32").Value = '  -2.76%  '
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("E34").Value = '  +5.78%  '
$ws.Range("D35").Value = '1.21'
$ws.Range("E35").Value = '  +6.37%  '
$ws.Range("E36").Value = '  +11.02%  '
$ws.Range("E37").Value = '  +1.71%  '
$ws.Range("D38").Value = '91.29'
$ws.Range("E38").Value = '  -1.67%  '
$ws.Range("D39").Value = '1.06'
$ws.Range("E39").Value = '  +5.47%  '
$ws.Range("D40").Value = '1.339.88'
$ws.Range("E40").Value = '  +2.27%  '
$ws.Range("E41").Value = '  -0.59%  '
$ws.Range("D42").Value = '14.58'
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  -1.46%  '
$ws.Range("E44").Value = '  -2.66%  '
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("E47").Value = '  +1.97%  '
$ws.Range("D48").Value = '2.017.92'
$ws.Range("E48").Value = '  +1.25%  '
$ws.Range("E49").Value = '  +5.15%  '
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("D51").Value = '3.28'
$ws.Range("E51").Value = '  +15.10%  '
